$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheets ---
$ws1.Name = "Displacement Cal"
$ws2.Name = "Shunt Cal"

# --- Populate "Shunt Cal" (sheet2) header/notes block ---
$ws2.Range("A1").Value = "Load Cell - Shunt Calibration"
$ws2.Range("A2").Value = "11/23/2010 - Large Scale LBCB 3"
$ws2.Range("A3").Value = "Ray, Michael, Alan, Ken, Party in the Crane Bay"

# --- Section header cells, mirroring the section titles used on "Displacement Cal" ---
# Same visual style as sheet1!A8 (bold 14pt Arial on pink fill) for the first five headers.
$ws2.Range("A5").Value = "X1 Dir."
$ws2.Range("A10").Value = "X2 Dir."
$ws2.Range("A14").Value = "Y1 Dir."
$ws2.Range("A18").Value = "Z1 Dir."
$ws2.Range("A22").Value = "Z2 Dir."
# Last header uses the style of sheet1!A106 (bold 14pt Arial on blue fill).
$ws2.Range("A26").Value = "Z3 Dir."

$ws1.Range("A8").Copy()
$ws2.Range("A5").PasteSpecial(-4122)
$ws1.Range("A8").Copy()
$ws2.Range("A10").PasteSpecial(-4122)
$ws1.Range("A8").Copy()
$ws2.Range("A14").PasteSpecial(-4122)
$ws1.Range("A8").Copy()
$ws2.Range("A18").PasteSpecial(-4122)
$ws1.Range("A8").Copy()
$ws2.Range("A22").PasteSpecial(-4122)

$ws1.Range("A106").Copy()
$ws2.Range("A26").PasteSpecial(-4122)

# Row heights for the header rows (matches the 18pt rows used on "Displacement Cal")
$ws2.Rows.Item(5).RowHeight = 18
$ws2.Rows.Item(10).RowHeight = 18
$ws2.Rows.Item(14).RowHeight = 18
$ws2.Rows.Item(18).RowHeight = 18
$ws2.Rows.Item(22).RowHeight = 18
$ws2.Rows.Item(26).RowHeight = 18

# --- View / selection state ---
$ws1.Activate()
$ws1.Range("A8:A122").Select()

$ws2.Activate()
$ws2.Range("A6").Select()
